$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

# Copy the date formatting (style) from an existing date cell (A2) to the
# new cell so the new row's Date column matches the existing rows exactly,
# then set the actual date value.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item($row, 1).Value = 42602.582511574074

$ws.Cells.Item($row, 2).Value = "Bag"
$ws.Cells.Item($row, 3).Value = 6203
$ws.Cells.Item($row, 4).Value = 9316
$ws.Cells.Item($row, 5).Value = 1145
$ws.Cells.Item($row, 6).Value = 133
$ws.Cells.Item($row, 7).Value = 61
$ws.Cells.Item($row, 8).Value = 67
$ws.Cells.Item($row, 9).Value = 31
$ws.Cells.Item($row, 10).Value = 3
$ws.Cells.Item($row, 11).Value = 3
$ws.Cells.Item($row, 12).Value = 49
$ws.Cells.Item($row, 13).Value = 49
